# Auto-generated Excel COM-interop edit script
# Applies scheduled-runner price/profit updates to the Leve profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 649.75
$ws.Range("I29").Value = 366.33334
$ws.Range("K29").Value = 1099.00002
$ws.Range("M29").Value = -818.0000199999999
$ws.Range("H40").Value = 2080
$ws.Range("H107").Value = 100
$ws.Range("I107").Value = 95
$ws.Range("J107").Value = 105
$ws.Range("K107").Value = 95
$ws.Range("L107").Value = 105
$ws.Range("M107").Value = 1825
$ws.Range("N107").Value = -3945
$ws.Range("H116").Value = 1551.25
$ws.Range("I116").Value = 1401.6666
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 1401.6666
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 2040.3334
$ws.Range("N116").Value = -8884
$ws.Range("H132").Value = 5839.3887
$ws.Range("I132").Value = 5928.067
$ws.Range("J132").Value = 5396
$ws.Range("K132").Value = 17784.201
$ws.Range("L132").Value = 16188
$ws.Range("M132").Value = -15254.201
$ws.Range("N132").Value = -21248
$ws.Range("H137").Value = 1947.5
$ws.Range("I137").Value = 1947.5
$ws.Range("K137").Value = 5842.5
$ws.Range("M137").Value = -3292.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 56
$ws.Range("I3").Value = 62
$ws.Range("K3").Value = 62
$ws.Range("M3").Value = 53
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = $null
$ws.Range("H13").Value = 99
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = $null
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = $null
$ws.Range("H74").Value = 2747.7144
$ws.Range("I74").Value = 2372.3333
$ws.Range("K74").Value = 2372.3333
$ws.Range("M74").Value = -1498.3333
$ws.Range("H77").Value = 2747.7144
$ws.Range("I77").Value = 2372.3333
$ws.Range("K77").Value = 11861.6665
$ws.Range("M77").Value = -7493.666499999999
$ws.Range("H97").Value = 1020.7143
$ws.Range("I97").Value = 899.1667
$ws.Range("K97").Value = 899.1667
$ws.Range("M97").Value = -403.1667
$ws.Range("H101").Value = 45001
$ws.Range("J101").Value = 45001
$ws.Range("L101").Value = 45001
$ws.Range("N101").Value = -51491
$ws.Range("H110").Value = 651.625
$ws.Range("I110").Value = 681.8570999999999
$ws.Range("J110").Value = 440
$ws.Range("K110").Value = 681.8570999999999
$ws.Range("L110").Value = 440
$ws.Range("M110").Value = 1363.1429
$ws.Range("N110").Value = -4530
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 99500180
$ws.Range("I7").Value = 99500180
$ws.Range("K7").Value = 99500180
$ws.Range("M7").Value = -99500067
$ws.Range("H86").Value = 3839.6667
$ws.Range("I86").Value = 3839.6667
$ws.Range("K86").Value = 3839.6667
$ws.Range("M86").Value = -2716.6667
$ws.Range("H89").Value = 3839.6667
$ws.Range("I89").Value = 3839.6667
$ws.Range("K89").Value = 19198.3335
$ws.Range("M89").Value = -13582.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 10000138
$ws.Range("I6").Value = 13333433
$ws.Range("K6").Value = 13333433
$ws.Range("M6").Value = -13333320
$ws.Range("H7").Value = 157.25
$ws.Range("I7").Value = 203
$ws.Range("K7").Value = 203
$ws.Range("M7").Value = -90
$ws.Range("H19").Value = 2643.1333
$ws.Range("I19").Value = 204.14285
$ws.Range("J19").Value = 4777.25
$ws.Range("K19").Value = 204.14285
$ws.Range("L19").Value = 4777.25
$ws.Range("M19").Value = -34.14285000000001
$ws.Range("N19").Value = -5117.25
$ws.Range("H24").Value = 2643.1333
$ws.Range("I24").Value = 204.14285
$ws.Range("J24").Value = 4777.25
$ws.Range("K24").Value = 204.14285
$ws.Range("L24").Value = 4777.25
$ws.Range("M24").Value = -34.14285000000001
$ws.Range("N24").Value = -5117.25
$ws.Range("H25").Value = 8897.125
$ws.Range("I25").Value = 8897.125
$ws.Range("K25").Value = 8897.125
$ws.Range("M25").Value = -8723.125
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = $null
$ws.Range("H107").Value = 709.4
$ws.Range("I107").Value = 709.8
$ws.Range("K107").Value = 709.8
$ws.Range("M107").Value = 1210.2
$ws.Range("H134").Value = 1536.2
$ws.Range("I134").Value = 1536.2
$ws.Range("K134").Value = 4608.6
$ws.Range("M134").Value = -2073.6
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 83
$ws.Range("I2").Value = 102.333336
$ws.Range("J2").Value = 63.666668
$ws.Range("K2").Value = 614.000016
$ws.Range("L2").Value = 382.000008
$ws.Range("M2").Value = -501.000016
$ws.Range("N2").Value = -608.000008
$ws.Range("H12").Value = 176
$ws.Range("I12").Value = 110.8
$ws.Range("K12").Value = 332.4
$ws.Range("M12").Value = -159.4
$ws.Range("H107").Value = 2241.2
$ws.Range("J107").Value = 2334.1667
$ws.Range("L107").Value = 7002.500100000001
$ws.Range("N107").Value = -10842.5001
$ws.Range("H121").Value = 1411.55
$ws.Range("J121").Value = 1419.4706
$ws.Range("L121").Value = 4258.4118
$ws.Range("N121").Value = -6878.4118
$ws.Range("H131").Value = 2826.4666
$ws.Range("J131").Value = 2826.4666
$ws.Range("L131").Value = 8479.399800000001
$ws.Range("N131").Value = -18559.3998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 24.5
$ws.Range("J2").Value = 25
$ws.Range("L2").Value = 25
$ws.Range("N2").Value = -251
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = $null
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").Value = $null
$ws.Range("H107").Value = 555
$ws.Range("I107").Value = 555
$ws.Range("K107").Value = 555
$ws.Range("M107").Value = 1365

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 7116.4
$ws.Range("I9").Value = 1860.6666
$ws.Range("J9").Value = 15000
$ws.Range("K9").Value = 1860.6666
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = -1636.6666
$ws.Range("N9").Value = -15448
$ws.Range("H46").Value = 1995.5
$ws.Range("I46").Value = 1995.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1995.5
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1807.5
$ws.Range("N46").Value = $null
$ws.Range("H61").Value = 1500
$ws.Range("I61").Value = 1500
$ws.Range("K61").Value = 1500
$ws.Range("M61").Value = -1298
$ws.Range("H82").Value = 1216.7142
$ws.Range("I82").Value = 1216.7142
$ws.Range("K82").Value = 1216.7142
$ws.Range("M82").Value = -855.7141999999999
$ws.Range("H85").Value = 1216.7142
$ws.Range("I85").Value = 1216.7142
$ws.Range("K85").Value = 1216.7142
$ws.Range("M85").Value = 31.28580000000011
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 65015
$ws.Range("J22").Value = 65015
$ws.Range("L22").Value = 65015
$ws.Range("N22").Value = -65601
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = $null
$ws.Range("H107").Value = 397.33334
$ws.Range("J107").Value = 550
$ws.Range("L107").Value = 1650
$ws.Range("N107").Value = -5490
$ws.Range("H136").Value = 4439.4287
$ws.Range("I136").Value = 4439.4287
$ws.Range("K136").Value = 13318.2861
$ws.Range("M136").Value = -10768.2861

